# Generate Report for Handoff
# A new localization hand-off ("8135d6b4-f305-45f0-b04b-dcbba50a3d30") is
# reported alongside the existing one ("87b58ab9-713d-404b-96d0-a38a3bbb52d0").
# A new row 2 is inserted on every sheet for the new file, pushing the
# previously-existing row down to row 3.

$wb = $excel.ActiveWorkbook

$oldUuid = "87b58ab9-713d-404b-96d0-a38a3bbb52d0"
$oldHash = "2993d00fe1daacb52c128f438d89ecd4cd7f0a7a"
$newUuid = "8135d6b4-f305-45f0-b04b-dcbba50a3d30"
$newHash = "08b003844bf439e9423fabc185d486aedb37af0b"

$oldMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/8cef19c92129b2b86acb926910ab90c269cc2b3f/e2e/$oldUuid.md"
$newMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/8cef19c92129b2b86acb926910ab90c269cc2b3f/e2e/$newUuid.md"

$oldZhCnXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/75bd3e72065f89080e838488d7a0140e4630ec1b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldUuid.$oldHash.zh-cn.xlf"
$newZhCnXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/75bd3e72065f89080e838488d7a0140e4630ec1b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newUuid.$newHash.zh-cn.xlf"

$oldDeDeXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5f15f84985fcf1be7775f48ac1203615032dc875/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldUuid.$oldHash.de-de.xlf"
$newDeDeXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5f15f84985fcf1be7775f48ac1203615032dc875/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newUuid.$newHash.de-de.xlf"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Sheets.Item("Overview")

$wsOverview.Hyperlinks.Delete()
$wsOverview.Rows.Item(2).Insert()

$wsOverview.Range("A2").Value = "$newUuid.md"
$wsOverview.Range("B2").Value = "zh-cn"
$wsOverview.Range("C2").Value = "de-de"
$wsOverview.Range("D2").Value = "2016-24-11 14:24:59"

$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-24-11 14:24:45"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $newMdUrl, "", "", "$newUuid.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $oldMdUrl, "", "", "$oldUuid.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Sheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Rows.Item(2).Insert()

$wsZhCn.Range("A2").Value = "$newUuid.md"
$wsZhCn.Range("B2").Value = ".md"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("D2").Value = "$newUuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("E2").Value = "2016-03-11 14:24:56"
$wsZhCn.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("H2").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("I2").Value = "Include"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $newMdUrl, "", "", "$newUuid.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B2"), $newMdUrl, "", "", ".md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), $newZhCnXlfUrl, "", "", "$newUuid.$newHash.zh-cn.xlf") | Out-Null

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $oldMdUrl, "", "", "$oldUuid.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B3"), $oldMdUrl, "", "", ".md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D3"), $oldZhCnXlfUrl, "", "", "$oldUuid.$oldHash.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Sheets.Item("de-de")

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Rows.Item(2).Insert()

$wsDeDe.Range("A2").Value = "$newUuid.md"
$wsDeDe.Range("B2").Value = ".md"
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("D2").Value = "$newUuid.$newHash.de-de.xlf"
$wsDeDe.Range("E2").Value = "2016-03-11 14:24:59"
$wsDeDe.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("H2").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("I2").Value = "Include"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $newMdUrl, "", "", "$newUuid.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B2"), $newMdUrl, "", "", ".md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), $newDeDeXlfUrl, "", "", "$newUuid.$newHash.de-de.xlf") | Out-Null

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $oldMdUrl, "", "", "$oldUuid.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B3"), $oldMdUrl, "", "", ".md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D3"), $oldDeDeXlfUrl, "", "", "$oldUuid.$oldHash.de-de.xlf") | Out-Null

Write-Output "Report generated for handoff: $newUuid"
